$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 255
$ws.Range("I2").Value = 768
$ws.Range("J2").Value = 3150
$ws.Range("K2").Value = 13
$ws.Range("L2").Value = 792
$ws.Range("M2").Value = 64
$ws.Range("N2").Value = 549
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 18
$ws.Range("Q2").Value = 8
$ws.Range("R2").Value = 44
$ws.Range("S2").Value = 346
$ws.Range("T2").Value = 557
$ws.Range("U2").Value = 34
$ws.Range("V2").Value = 4773
$ws.Range("X2").Value = 4708
$ws.Range("Y2").Value = 7
$ws.Range("Z2").Value = 57
$ws.Range("AA2").Value = 20
